# Updated working hours: add a new entry for 17.9.2025 (row 9)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Date for the new row
$ws.Range("A9").Value = "17.9.2025"

# Start/End time pairs (stored as Excel time fractions, displayed with existing time format)
$ws.Range("B9").Value = 0.58333333333333337
$ws.Range("C9").Value = 0.64583333333333337
$ws.Range("D9").Value = 0.89583333333333337
$ws.Range("E9").Value = 0.96527777777777779

# Match the time number format used by the rest of the sheet (column B/C/D/E of existing rows)
$ws.Range("B9:E9").NumberFormat = $ws.Range("B8:E8").NumberFormat

# Move active selection to F12 as in the saved workbook
$ws.Range("F12").Select()
